$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.897.88"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.632.55"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "211.37"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "23.46"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").Value = "0.0882"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "1.864.25"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "1.634.04"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "0.565"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "65.44"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "27.909.26"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "228.95"
$ws.Range("D19").Value = "7.68"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "10.06"
$ws.Range("E23").Value = "  -3.35%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "155.29"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "15.53"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").Value = "3.12"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("D34").Value = "1.392.04"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("E36").Value = "  +9.56%  "
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").Value = "0.558"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").Value = "65.76"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").Value = "5.43"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("D46").Value = "1.773.61"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("D48").Value = "88.78"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").Value = "0.0504"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").Value = "  +0.69%  "
